$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ B=...; C=...; D=...; E=... } (only changed columns present)
$updates = @{
    2 = @{ D="29.439.54"; E="  +0.13%  " }
    3 = @{ D="1.893.01"; E="  -1.20%  " }
    4 = @{ D="1.003"; E="  -0.61%  " }
    5 = @{ D="338.37"; E="  +4.21%  " }
    6 = @{ D="1.003"; E="  -0.41%  " }
    7 = @{ D="0.4751"; E="  -1.28%  " }
    8 = @{ D="0.3997"; E="  -1.48%  " }
    9 = @{ D="47.03"; E="  -1.95%  " }
    10 = @{ D="0.08016"; E="  -2.34%  " }
    11 = @{ D="0.9882"; E="  -1.86%  " }
    12 = @{ D="23.12"; E="  -0.52%  " }
    13 = @{ D="1.922.96"; E="  +0.85%  " }
    14 = @{ D="5.921"; E="  -2.35%  " }
    15 = @{ D="7.071"; E="  -2.03%  " }
    16 = @{ D="88.96"; E="  -2.86%  " }
    17 = @{ D="0.06803"; E="  -0.94%  " }
    18 = @{ E="  -0.23%  " }
    19 = @{ D="0.00001018"; E="  -1.92%  " }
    20 = @{ D="17.29"; E="  -1.66%  " }
    21 = @{ D="1.002"; E="  -0.52%  " }
    22 = @{ D="29.478.96"; E="  +0.22%  " }
    23 = @{ D="5.506"; E="  -2.61%  " }
    24 = @{ D="11.63"; E="  -1.02%  " }
    25 = @{ D="2.150"; E="  -1.59%  " }
    26 = @{ D="2.155.95"; E="  +0.86%  " }
    27 = @{ D="157.27"; E="  +1.00%  " }
    28 = @{ D="6.469"; E="  -1.41%  " }
    29 = @{ D="19.61"; E="  -1.50%  " }
    30 = @{ D="2.046"; E="  -3.02%  " }
    31 = @{ D="118.85"; E="  -1.36%  " }
    32 = @{ D="0.9941"; E="  -1.86%  " }
    33 = @{ D="0.09539"; E="  -0.88%  " }
    34 = @{ D="5.466"; E="  -2.83%  " }
    35 = @{ D="1.383"; E="  +1.02%  " }
    36 = @{ D="3.527"; E="  -0.57%  " }
    37 = @{ D="0.06414"; E="  +5.29%  " }
    38 = @{ D="0.02241"; E="  -1.84%  " }
    39 = @{ D="1.200"; E="  +1.56%  " }
    40 = @{ D="0.5832"; E="  -1.85%  " }
    41 = @{ D="10.52"; E="  -3.47%  " }
    42 = @{ D="7.742"; E="  -3.65%  " }
    43 = @{ D="0.1818"; E="  -1.48%  " }
    44 = @{ D="2.427"; E="  +1.94%  " }
    45 = @{ D="1.264"; E="  -1.18%  " }
    46 = @{ D="12.12"; E="  -2.67%  " }
    47 = @{ D="0.5491"; E="  -1.52%  " }
    48 = @{ D="0.07333"; E="  -3.55%  " }
    49 = @{ D="1.949"; E="  +0.16%  " }
    50 = @{ D="116.33"; E="  -1.72%  " }
    51 = @{ B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.369"; E="  -2.42%  " }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    foreach ($col in $cells.Keys) {
        $addr = "$col$row"
        $value = $cells[$col]
        if ($col -eq "D") {
            # Force text storage so numeric-looking strings (e.g. "1.003",
            # "0.08016") are not reinterpreted as numbers by Excel.
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value = $value
            $ws.Range($addr).ClearFormats()
        } else {
            $ws.Range($addr).Value = $value
        }
    }
}
